$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Client" column: it moves from E to after "Typist QC" (new column G) ---
# Cutting column E and inserting it (with shift) right before the (then) column H
# shifts Typist/Typist QC left into E/F and drops Client into G, carrying the
# original cell formatting (incl. the highlighted style) and column width along.
$ws.Columns("E").Cut()
$ws.Columns("H").Insert(-4161)

# --- Swap "Product Name" (now H) and "Lob" (now J) while leaving "Process" (I) in place ---
$ws.Range("H1:H3").Copy()
$ws.Range("Z1:Z3").PasteSpecial(-4104)
$ws.Range("J1:J3").Copy()
$ws.Range("H1:H3").PasteSpecial(-4104)
$ws.Range("Z1:Z3").Copy()
$ws.Range("J1:J3").PasteSpecial(-4104)
$ws.Range("Z1:Z3").Clear()

# --- Re-home the custom column width that used to live on column E onto column G ---
# (The cut/insert above leaves stray zero-width leftovers on E:F; clear those out and
#  restore the plain formatting that belongs on the header/data rows there.)
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)

$ws.Columns("E:F").ClearFormats()

$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)

# --- Update the active cell selection on the sheet ---
$ws.Range("F4").Select()
